$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "o554F"

# Add new row 16 data (mirrors the pattern of previous rows, HexGrid-60degTilt5degRes entry)
$ws.Range("A16").Value = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16:M16").Value = 1
